$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the subscript-three character (U+2083) used in cell D19 safely,
# avoiding PowerShell's numeric coercion when concatenating with "+".
$sub3 = [char]0x2083

$ws.Range("D2").Value = "27.549.17"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.647.89"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.69"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.534"
$ws.Range("E6").Value = "  +5.15%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.57"
$ws.Range("E8").Value = "  -2.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.257"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0611"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0892"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "1.881.65"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").Value = "1.644.49"
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("E14").Value = "  +4.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.04"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.56"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").Value = "27.533.56"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.17"
$ws.Range("E18").Value = "  -3.89%  "
$ws.Range("D19").Value = ("0.0{0}0724" -f $sub3)
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.32"
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.75"
$ws.Range("E23").Value = "  +4.06%  "
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.84"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.02"
$ws.Range("E26").Value = "  -2.43%  "
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.61"
$ws.Range("E29").Value = "  -4.19%  "
$ws.Range("E30").Value = "  -2.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0487"
$ws.Range("E31").Value = "  -3.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.20"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("D34").Value = "1.430.27"
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.568"
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.885"
$ws.Range("E38").Value = "  -4.31%  "
$ws.Range("E39").Value = "  -2.86%  "
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.818"
$ws.Range("E42").Value = "  +3.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.56"
$ws.Range("E43").Value = "  +2.60%  "
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.27"
$ws.Range("E45").Value = "  -6.44%  "
$ws.Range("D46").Value = "1.791.12"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.69"
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.09"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.79"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0991"
$ws.Range("E51").Value = "  -3.36%  "
